# 9.c.1 update: refresh footnote source (MCR KR instead of GKITS KR), add 2023 column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the 2023 data column (O) -------------------------------------------------
# Copy formatting from the 2022 column (N) for the header band (row3), the year
# header (row4) and the three data rows (5-7), then overwrite with the 2023 values.
$ws.Range("N3:N7").Copy($ws.Range("O3:O7"))

$ws.Range("O4").Value = 2023
$ws.Range("O5").Value = 99
$ws.Range("O6").Value = 98.9
$ws.Range("O7").Value = 98.8

# Widen the label columns a touch now that the table is wider.
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 37.166666666666664

# --- Refresh the footnote text (new source agency) ---------------------------------
$ws.Range("B8").Value = "*по данным МЦР КР"
$ws.Range("C8").Value = "*according to the MDD KR"
$ws.Range("A8").Value = "*КР СӨМ маалыматтары  боюнча"

# Leave the selection parked on A1 (matches the saved/default view state).
$ws.Range("A1").Select() | Out-Null
